$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the header cell A1 from "RUN_TEST" to "RUN"
$ws.Range("A1").Value = "RUN"

# Add the missing row 5, duplicating the pattern of row 2/3 (RUN / student / Password123)
$ws.Range("A5").Value = "RUN"
$ws.Range("B5").Value = "student"
$ws.Range("C5").Value = "Password123"

# Update the active selection to A4 (matches the post-edit selection in the file)
$ws.Range("A4").Select()
